$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update regression coefficients for the 12-lag monthly diffs/rates tests.
# Row 2 = "A Lag", Row 3 = "C Lag"; Column B = "A", Column C = "C"
$ws.Range("B2").Value = "-0.372***"
$ws.Range("B3").Value = "-3.464***"
$ws.Range("C2").Value = "0.01*"
$ws.Range("C3").Value = "-0.808***"
